$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(46, 8).Value = 169999
$ws.Cells.Item(46, 10).Value = 169999
$ws.Cells.Item(46, 12).Value = 509997
$ws.Cells.Item(46, 14).Value = -510235
$ws.Cells.Item(60, 8).Value = 169999
$ws.Cells.Item(60, 10).Value = 169999
$ws.Cells.Item(60, 12).Value = 509997
$ws.Cells.Item(60, 14).Value = -510965
$ws.Cells.Item(113, 8).Value = 10123.6
$ws.Cells.Item(113, 9).Value = 8433.833000000001
$ws.Cells.Item(113, 10).Value = 12658.25
$ws.Cells.Item(113, 11).Value = 8433.833000000001
$ws.Cells.Item(113, 12).Value = 12658.25
$ws.Cells.Item(113, 13).Value = -5179.833000000001
$ws.Cells.Item(113, 14).Value = -19166.25
$ws.Cells.Item(125, 8).Value = 7421.222
$ws.Cells.Item(125, 9).Value = 17998
$ws.Cells.Item(125, 11).Value = 161982
$ws.Cells.Item(125, 13).Value = -159522
$ws.Cells.Item(137, 8).Value = 9856.385
$ws.Cells.Item(137, 9).Value = 17548.691
$ws.Cells.Item(137, 10).Value = 2164.077
$ws.Cells.Item(137, 11).Value = 52646.073
$ws.Cells.Item(137, 12).Value = 6492.231000000001
$ws.Cells.Item(137, 13).Value = -50096.073
$ws.Cells.Item(137, 14).Value = -11592.231
$ws.Cells.Item(138, 8).Value = 225791.19
$ws.Cells.Item(138, 9).Value = 507308.22
$ws.Cells.Item(138, 10).Value = 3989.879
$ws.Cells.Item(138, 11).Value = 1521924.66
$ws.Cells.Item(138, 12).Value = 11969.637
$ws.Cells.Item(138, 13).Value = -1516784.66
$ws.Cells.Item(138, 14).Value = -22249.637

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 7490.8335
$ws.Cells.Item(31, 9).Value = 7490.8335
$ws.Cells.Item(31, 11).Value = 7490.8335
$ws.Cells.Item(31, 13).Value = -7196.8335
$ws.Cells.Item(45, 8).Value = 52366.926
$ws.Cells.Item(45, 9).Value = 77198.41
$ws.Cells.Item(45, 10).Value = 4477.643
$ws.Cells.Item(45, 11).Value = 77198.41
$ws.Cells.Item(45, 12).Value = 4477.643
$ws.Cells.Item(45, 13).Value = -76821.41
$ws.Cells.Item(45, 14).Value = -5231.643
$ws.Cells.Item(74, 8).Value = 6496
$ws.Cells.Item(74, 9).Value = 7788
$ws.Cells.Item(74, 10).Value = 3542.8572
$ws.Cells.Item(74, 11).Value = 7788
$ws.Cells.Item(74, 12).Value = 3542.8572
$ws.Cells.Item(74, 13).Value = -6914
$ws.Cells.Item(74, 14).Value = -5290.8572
$ws.Cells.Item(77, 8).Value = 6496
$ws.Cells.Item(77, 9).Value = 7788
$ws.Cells.Item(77, 10).Value = 3542.8572
$ws.Cells.Item(77, 11).Value = 38940
$ws.Cells.Item(77, 12).Value = 17714.286
$ws.Cells.Item(77, 13).Value = -34572
$ws.Cells.Item(77, 14).Value = -26450.286
$ws.Cells.Item(88, 8).Value = 55556588
$ws.Cells.Item(88, 9).Value = 537
$ws.Cells.Item(88, 10).Value = 90910440
$ws.Cells.Item(88, 11).Value = 537
$ws.Cells.Item(88, 12).Value = 90910440
$ws.Cells.Item(88, 13).Value = -131
$ws.Cells.Item(88, 14).Value = -90911252
$ws.Cells.Item(91, 8).Value = 55556588
$ws.Cells.Item(91, 9).Value = 537
$ws.Cells.Item(91, 10).Value = 90910440
$ws.Cells.Item(91, 11).Value = 537
$ws.Cells.Item(91, 12).Value = 90910440
$ws.Cells.Item(91, 13).Value = 867
$ws.Cells.Item(91, 14).Value = -90913248
$ws.Cells.Item(110, 8).Value = 2230.238
$ws.Cells.Item(110, 10).Value = 4625
$ws.Cells.Item(110, 12).Value = 4625
$ws.Cells.Item(110, 14).Value = -8715

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3461.4211
$ws.Cells.Item(20, 9).Value = 2083.818
$ws.Cells.Item(20, 10).Value = 5355.625
$ws.Cells.Item(20, 11).Value = 2083.818
$ws.Cells.Item(20, 12).Value = 5355.625
$ws.Cells.Item(20, 13).Value = -1836.818
$ws.Cells.Item(20, 14).Value = -5849.625
$ws.Cells.Item(75, 8).Value = 85487.25
$ws.Cells.Item(75, 9).Value = 80649.664
$ws.Cells.Item(75, 11).Value = 80649.664
$ws.Cells.Item(75, 13).Value = -79713.664
$ws.Cells.Item(78, 8).Value = 85487.25
$ws.Cells.Item(78, 9).Value = 80649.664
$ws.Cells.Item(78, 11).Value = 241948.992
$ws.Cells.Item(78, 13).Value = -237268.992
$ws.Cells.Item(86, 8).Value = 5329.087
$ws.Cells.Item(86, 9).Value = 8360.333000000001
$ws.Cells.Item(86, 11).Value = 8360.333000000001
$ws.Cells.Item(86, 13).Value = -7237.333000000001
$ws.Cells.Item(89, 8).Value = 5329.087
$ws.Cells.Item(89, 9).Value = 8360.333000000001
$ws.Cells.Item(89, 11).Value = 41801.665
$ws.Cells.Item(89, 13).Value = -36185.665
$ws.Cells.Item(134, 8).Value = 7076.174
$ws.Cells.Item(134, 9).Value = 7664
$ws.Cells.Item(134, 10).Value = 3157.3333
$ws.Cells.Item(134, 11).Value = 22992
$ws.Cells.Item(134, 12).Value = 9471.999899999999
$ws.Cells.Item(134, 13).Value = -20457
$ws.Cells.Item(134, 14).Value = -14541.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10486.294
$ws.Cells.Item(31, 9).Value = 11945.615
$ws.Cells.Item(31, 10).Value = 5743.5
$ws.Cells.Item(31, 11).Value = 11945.615
$ws.Cells.Item(31, 12).Value = 5743.5
$ws.Cells.Item(31, 13).Value = -11650.615
$ws.Cells.Item(31, 14).Value = -6333.5
$ws.Cells.Item(34, 8).Value = 10486.294
$ws.Cells.Item(34, 9).Value = 11945.615
$ws.Cells.Item(34, 10).Value = 5743.5
$ws.Cells.Item(34, 11).Value = 11945.615
$ws.Cells.Item(34, 12).Value = 5743.5
$ws.Cells.Item(34, 13).Value = -11743.615
$ws.Cells.Item(34, 14).Value = -6147.5
$ws.Cells.Item(122, 8).Value = 10794.846
$ws.Cells.Item(122, 9).Value = 16282.875
$ws.Cells.Item(122, 11).Value = 48848.625
$ws.Cells.Item(122, 13).Value = -46398.625
$ws.Cells.Item(134, 8).Value = 7823.524
$ws.Cells.Item(134, 9).Value = 9986.4
$ws.Cells.Item(134, 11).Value = 29959.2
$ws.Cells.Item(134, 13).Value = -27424.2
$ws.Cells.Item(141, 8).Value = 311421.94
$ws.Cells.Item(141, 10).Value = 400121.1
$ws.Cells.Item(141, 12).Value = 400121.1
$ws.Cells.Item(141, 14).Value = -410481.1

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 8731.362999999999
$ws.Cells.Item(81, 10).Value = 8731.362999999999
$ws.Cells.Item(81, 12).Value = 26194.089
$ws.Cells.Item(81, 14).Value = -28440.089
$ws.Cells.Item(84, 8).Value = 8731.362999999999
$ws.Cells.Item(84, 10).Value = 8731.362999999999
$ws.Cells.Item(84, 12).Value = 78582.26699999999
$ws.Cells.Item(84, 14).Value = -89814.26699999999
$ws.Cells.Item(103, 8).Value = 5366.5835
$ws.Cells.Item(103, 9).Value = 10411.8
$ws.Cells.Item(103, 10).Value = 1762.8572
$ws.Cells.Item(103, 11).Value = 31235.4
$ws.Cells.Item(103, 12).Value = 5288.571599999999
$ws.Cells.Item(103, 13).Value = -30356.4
$ws.Cells.Item(103, 14).Value = -7046.571599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 36330
$ws.Cells.Item(63, 9).Value = 29495
$ws.Cells.Item(63, 10).Value = 50000
$ws.Cells.Item(63, 11).Value = 29495
$ws.Cells.Item(63, 12).Value = 50000
$ws.Cells.Item(63, 13).Value = -28809
$ws.Cells.Item(63, 14).Value = -51372
$ws.Cells.Item(66, 8).Value = 36330
$ws.Cells.Item(66, 9).Value = 29495
$ws.Cells.Item(66, 10).Value = 50000
$ws.Cells.Item(66, 11).Value = 88485
$ws.Cells.Item(66, 12).Value = 150000
$ws.Cells.Item(66, 13).Value = -85053
$ws.Cells.Item(66, 14).Value = -156864
$ws.Cells.Item(70, 8).Value = 7662.826
$ws.Cells.Item(70, 9).Value = 6549.4707
$ws.Cells.Item(70, 11).Value = 6549.4707
$ws.Cells.Item(70, 13).Value = -6279.4707
$ws.Cells.Item(73, 8).Value = 7662.826
$ws.Cells.Item(73, 9).Value = 6549.4707
$ws.Cells.Item(73, 11).Value = 6549.4707
$ws.Cells.Item(73, 13).Value = -5613.4707
$ws.Cells.Item(102, 8).Value = 6432.769
$ws.Cells.Item(102, 9).Value = 7313.7144
$ws.Cells.Item(102, 10).Value = 2732.8
$ws.Cells.Item(102, 11).Value = 7313.7144
$ws.Cells.Item(102, 12).Value = 2732.8
$ws.Cells.Item(102, 13).Value = -5691.7144
$ws.Cells.Item(102, 14).Value = -5976.8
$ws.Cells.Item(122, 8).Value = 12148.944
$ws.Cells.Item(122, 9).Value = 8343.546
$ws.Cells.Item(122, 10).Value = 18128.857
$ws.Cells.Item(122, 11).Value = 25030.638
$ws.Cells.Item(122, 12).Value = 54386.571
$ws.Cells.Item(122, 13).Value = -22580.638
$ws.Cells.Item(122, 14).Value = -59286.571
$ws.Cells.Item(132, 8).Value = 3611.22
$ws.Cells.Item(132, 9).Value = 3736.3416
$ws.Cells.Item(132, 11).Value = 11209.0248
$ws.Cells.Item(132, 13).Value = -8679.024800000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 19574.188
$ws.Cells.Item(40, 9).Value = 23403.15
$ws.Cells.Item(40, 10).Value = 13192.583
$ws.Cells.Item(40, 11).Value = 23403.15
$ws.Cells.Item(40, 12).Value = 13192.583
$ws.Cells.Item(40, 13).Value = -23267.15
$ws.Cells.Item(40, 14).Value = -13464.583
$ws.Cells.Item(122, 8).Value = 4704.0293
$ws.Cells.Item(122, 9).Value = 4529.72
$ws.Cells.Item(122, 10).Value = 5188.222
$ws.Cells.Item(122, 11).Value = 13589.16
$ws.Cells.Item(122, 12).Value = 15564.666
$ws.Cells.Item(122, 13).Value = -11139.16
$ws.Cells.Item(122, 14).Value = -20464.666
$ws.Cells.Item(132, 8).Value = 1067377.9
$ws.Cells.Item(132, 9).Value = 1657134.8
$ws.Cells.Item(132, 10).Value = 5815.6
$ws.Cells.Item(132, 11).Value = 4971404.4
$ws.Cells.Item(132, 12).Value = 17446.8
$ws.Cells.Item(132, 13).Value = -4968874.4
$ws.Cells.Item(132, 14).Value = -22506.8
$ws.Cells.Item(136, 8).Value = 3892.2
$ws.Cells.Item(136, 9).Value = 2855.6155
$ws.Cells.Item(136, 11).Value = 8566.8465
$ws.Cells.Item(136, 13).Value = -6016.8465

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 90925240
$ws.Cells.Item(2, 9).Value = 111119960
$ws.Cells.Item(2, 11).Value = 111119960
$ws.Cells.Item(2, 13).Value = -111119848
$ws.Cells.Item(4, 8).Value = 1077.2727
$ws.Cells.Item(4, 9).Value = 407.14285
$ws.Cells.Item(4, 11).Value = 407.14285
$ws.Cells.Item(4, 13).Value = -294.14285
$ws.Cells.Item(107, 8).Value = 18988.766
$ws.Cells.Item(107, 10).Value = 100094.336
$ws.Cells.Item(107, 12).Value = 300283.008
$ws.Cells.Item(107, 14).Value = -304123.008
$ws.Cells.Item(126, 8).Value = 37160.918
$ws.Cells.Item(126, 10).Value = 6326.1665
$ws.Cells.Item(126, 12).Value = 18978.4995
$ws.Cells.Item(126, 14).Value = -23918.4995
$ws.Cells.Item(136, 8).Value = 436787
$ws.Cells.Item(136, 10).Value = 21555.5
$ws.Cells.Item(136, 12).Value = 64666.5
$ws.Cells.Item(136, 14).Value = -69766.5
